$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Human Resources  Position Catagories Add Functionality", "FAILED", "chrome"),
    @("Human Resources  Position Catagories Add Functionality", "FAILED", "chrome"),
    @("Human Resources  Position Catagories Add Functionality", "FAILED", "chrome"),
    @("Human Resources  Position Catagories Add Functionality", "PASSED", "chrome"),
    @("Scholls  Position Catagories Add Functionality", "FAILED", "chrome"),
    @("Schools  Position Catagories Edit Functionality", "FAILED", "chrome"),
    @("Schools  Position Catagories Edit Functionality", "PASSED", "chrome"),
    @("Schools  Position Catagories Edit Functionality", "FAILED", "chrome"),
    @("Schools  Position Catagories Edit Functionality", "FAILED", "chrome"),
    @("Schools  Position Catagories Edit Functionality", "FAILED", "chrome")
)

$startRow = 24
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
